$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three new case rows appended below the existing table (rows 9-11).
# Column A (ID) is a true number; columns B-P mirror the sheet's existing
# convention of storing every other value (including numeric-looking docket
# numbers and ISO-style dates) as literal text.
$newRows = @(
    @{ Row=9;  A=14; B='3232';    C='2026-01-07'; D='keynn';  E='trekk'; F='San Isidro';    G='robbery';      H='Invalid Date'; I='2026-01-08'; J='Atty.batman';     K='Reppublic Act 2104'; L='Tagbi';   M='2026-01-09'; N='pending';    O='14k';  P='N/A' },
    @{ Row=10; A=15; B='245232';  C='2026-01-07'; D='koynnn'; E='kart';  F='Tiptip ';       G='grave threat'; H='Invalid Date'; I='2026-01-30'; J='Atty. Superman'; K='N/A';                L='tagbi';   M='2026-01-29'; N='terminated'; O='12k';  P='N/A' },
    @{ Row=11; A=16; B='4121214'; C='2026-01-29'; D='shaaan'; E='kroel'; F='Dauis';         G='raped';        H='Invalid Date'; I='2026-01-13'; J='Atty. Doroy';    K='N/A';                L='Capitol'; M='2026-01-22'; N='pending';    O='15k';  P='N/A' }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # A: numeric ID
    $ws.Cells.Item($row, 1).Value = $r.A

    # B-P: force text storage (leading apostrophe = "store as text"), then
    # strip the resulting quote-prefix cell format so no style is left
    # behind, matching the rest of the sheet.
    $textCols = @{ 2=$r.B; 3=$r.C; 4=$r.D; 5=$r.E; 6=$r.F; 7=$r.G; 8=$r.H; 9=$r.I; 10=$r.J; 11=$r.K; 12=$r.L; 13=$r.M; 14=$r.N; 15=$r.O; 16=$r.P }
    foreach ($col in $textCols.Keys) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = "'" + $textCols[$col]
        $cell.ClearFormats()
    }
}
